$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (content now matches former row 16)
$ws.Cells.Item(15, 1).Value = 111479776
$ws.Cells.Item(15, 2).Value = 96368
$ws.Cells.Item(15, 4).Value = "LC"
$ws.Cells.Item(15, 5).Value = 221952
$ws.Cells.Item(15, 6).Value = "Spindelblomster"
$ws.Cells.Item(15, 7).Value = "Neottia cordata"
$ws.Cells.Item(15, 8).Value = "(L.) Rich."
$ws.Cells.Item(15, 17).Value = 489790.1442766788
$ws.Cells.Item(15, 18).Value = 7140835.667408227
$ws.Cells.Item(15, 49).Value = "Jonathan Frendel"
$ws.Cells.Item(15, 50).Value = "Jonathan Frendel, Astrid Blomberg, Elias Blad, Elicia Olsson, Elvira Klang, Filippa Paperin, Iris Elmér, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg"

# Row 16 (content now matches former row 17)
$ws.Cells.Item(16, 1).Value = 111479784
$ws.Cells.Item(16, 2).Value = 78612
$ws.Cells.Item(16, 5).Value = 6464
$ws.Cells.Item(16, 6).Value = "Luddlav"
$ws.Cells.Item(16, 7).Value = "Nephroma resupinatum"
$ws.Cells.Item(16, 8).Value = "(L.) Ach."
$ws.Cells.Item(16, 17).Value = 489806.0819462601
$ws.Cells.Item(16, 18).Value = 7140678.121827234

# Row 17 (content now matches former row 15)
$ws.Cells.Item(17, 1).Value = 111476452
$ws.Cells.Item(17, 2).Value = 96348
$ws.Cells.Item(17, 4).Value = "VU"
$ws.Cells.Item(17, 5).Value = 220787
$ws.Cells.Item(17, 6).Value = "Knärot"
$ws.Cells.Item(17, 7).Value = "Goodyera repens"
$ws.Cells.Item(17, 8).Value = "(L.) R. Br."
$ws.Cells.Item(17, 17).Value = 489815.0024043967
$ws.Cells.Item(17, 18).Value = 7140755.108741223
$ws.Cells.Item(17, 49).Value = "Signe Propst"
$ws.Cells.Item(17, 50).Value = "Signe Propst"

# Row 18 (content now matches former row 19)
$ws.Cells.Item(18, 1).Value = 111476434
$ws.Cells.Item(18, 2).Value = 77515
$ws.Cells.Item(18, 5).Value = 6425
$ws.Cells.Item(18, 6).Value = "Garnlav"
$ws.Cells.Item(18, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(18, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(18, 17).Value = 490140.7149260837
$ws.Cells.Item(18, 18).Value = 7140869.131914983
$ws.Cells.Item(18, 19).Value = 10
$ws.Cells.Item(18, 49).Value = "Signe Propst"
$ws.Cells.Item(18, 50).Value = "Signe Propst"

# Row 19 (content now matches former row 20)
$ws.Cells.Item(19, 1).Value = 111476447
$ws.Cells.Item(19, 2).Value = 78578
$ws.Cells.Item(19, 5).Value = 6458
$ws.Cells.Item(19, 6).Value = "Lunglav"
$ws.Cells.Item(19, 7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(19, 8).Value = "(L.) Hoffm."
$ws.Cells.Item(19, 17).Value = 489901.6080066679
$ws.Cells.Item(19, 18).Value = 7140875.102526958

# Row 20 (content now matches former row 18)
$ws.Cells.Item(20, 1).Value = 111481809
$ws.Cells.Item(20, 17).Value = 489716.1870822187
$ws.Cells.Item(20, 18).Value = 7141000.754049809
$ws.Cells.Item(20, 19).Value = 15
$ws.Cells.Item(20, 49).Value = "Elias Blad"
$ws.Cells.Item(20, 50).Value = "Elias Blad, Astrid Blomberg, Elicia Olsson, Elvira Klang, Filippa Paperin, Iris Elmér, Ivar Anderberg, Jonathan Frendel, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg"

# Row 28 (content now matches former row 29)
$ws.Cells.Item(28, 1).Value = 111481799
$ws.Cells.Item(28, 2).Value = 90087
$ws.Cells.Item(28, 4).Value = "LC"
$ws.Cells.Item(28, 5).Value = 3298
$ws.Cells.Item(28, 6).Value = "Trådticka"
$ws.Cells.Item(28, 7).Value = "Climacocystis borealis"
$ws.Cells.Item(28, 8).Value = "(Fr.) Kotl. & Pouzar"
$ws.Cells.Item(28, 17).Value = 489701.4212471162
$ws.Cells.Item(28, 18).Value = 7140989.555053207

# Row 29 (content now matches former row 28)
$ws.Cells.Item(29, 1).Value = 111481804
$ws.Cells.Item(29, 2).Value = 78578
$ws.Cells.Item(29, 4).Value = "NT"
$ws.Cells.Item(29, 5).Value = 6458
$ws.Cells.Item(29, 6).Value = "Lunglav"
$ws.Cells.Item(29, 7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(29, 8).Value = "(L.) Hoffm."
$ws.Cells.Item(29, 17).Value = 489649.3234224396
$ws.Cells.Item(29, 18).Value = 7140827.054006468

